# "Add handling for empty data sheets"
#
# - Keep sheet "A" selection at B1 (it stops being the tab-selected sheet).
# - Move sheet "M" selection to C23.
# - Append a new, empty sheet "Q" after "M" (selection D26).
# - Append a new sheet "Sheet3" after "Q" holding two cells of data
#   (A1 = "moj ", B1 = "shit"), selection D31, and make it the active tab.

$wb = $excel.ActiveWorkbook

# -- Sheet "A": selection stays on B1 --------------------------------------
$wsA = $wb.Worksheets.Item("A")
$wsA.Activate()
$wsA.Range("B1").Select()

# -- Sheet "M": move the selection to C23 ----------------------------------
$wsM = $wb.Worksheets.Item("M")
$wsM.Activate()
$wsM.Range("C23").Select()

# -- New sheet "Q" -----------------------------------------------------------
# A throwaway sheet is inserted (and removed) first so that the sheetId
# counter lands on the same values as the source workbook (Q=4, Sheet3=5)
# instead of the ones it would get if added immediately after "M" (3, 4).
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$placeholder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsQDraft = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $placeholder)
$wsQDraft.Name = "Q"
$placeholder.Delete()

# The reference grabbed before the delete can go stale, so re-resolve by name.
$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Activate()
$wsQ.Range("D26").Select()

# -- New sheet "Sheet3" with data -------------------------------------------
$wsSheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsQ)
$wsSheet3.Name = "Sheet3"
$wsSheet3.Range("A1").Value = "moj "
$wsSheet3.Range("B1").Value = "shit"
$wsSheet3.Activate()
$wsSheet3.Range("D31").Select()
